$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the distinctive "last row" date format currently on row 50 (A50)
# so the new last row (row 51) can inherit it.
$lastRowDateFormat = $ws.Range("A50").NumberFormat

# Row 50 is no longer the last row, so give it the same date format used by
# all the other preceding data rows (e.g. row 49).
$ws.Range("A50").NumberFormat = $ws.Range("A49").NumberFormat

# Append the new day's results as row 51, the new last row.
$ws.Range("A51").Value = 45636
$ws.Range("B51").Value = 128
$ws.Range("C51").Value = 114
$ws.Range("D51").Value = 119

# Give the new last row the distinctive date format that row 50 used to have.
$ws.Range("A51").NumberFormat = $lastRowDateFormat
